# additional-units-conventional.xlsx - "Capacity" sheet update
# - updated electrolysis / hydrogen processor capacity figures
# - updated EV smart charger capacity for FI00
# - changed the Node AutoFilter from SE03 to FI00 (which re-derives the
#   hidden/visible state of every data row under the filtered column)
# - moved the selection/active cell to H94

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated capacity values (Other_capa column, H) ---
$ws.Range("H7").Value = 5000     # FI00 / Electrolysis
$ws.Range("H9").Value = 3500     # FI00 / Hydrogen processor
$ws.Range("H93").Value = 57200   # FI00 / EV smart charger

# --- Re-apply the Node AutoFilter so it shows FI00 instead of SE03 ---
# This also recalculates which rows are hidden/visible for every other
# Node value (e.g. rows for SE03 become hidden, rows for FI00 become shown).
[void]$ws.Range("A1:J99").AutoFilter(1, "FI00")

# --- Selection / active cell ---
[void]$ws.Range("H94").Select()
